$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 18.06930079656841
$ws.Range("C2").Value = 9.353341632968561
$ws.Range("D2").Value = 6.007124275456278
$ws.Range("E2").Value = 10.47931659354342
$ws.Range("G2").Value = 44.73420541563187
$ws.Range("H2").Value = 17.8721191610929
$ws.Range("I2").Value = 27.08284930237485
$ws.Range("L2").Value = 10.19120048263675
$ws.Range("M2").Value = 16.24703669116382
$ws.Range("N2").Value = 19.21324923203129
$ws.Range("B3").Value = 17.63239344643055
$ws.Range("C3").Value = 8.794550449536404
$ws.Range("D3").Value = 5.893121094777782
$ws.Range("E3").Value = 10.47807820660055
$ws.Range("G3").Value = 44.42195685623381
$ws.Range("H3").Value = 17.88449843640579
$ws.Range("I3").Value = 27.12715772226077
$ws.Range("L3").Value = 10.20135640171817
$ws.Range("M3").Value = 16.16692492213395
$ws.Range("N3").Value = 19.28554343610849
$ws.Range("B4").Value = 17.36322426951491
$ws.Range("C4").Value = 8.431349878735141
$ws.Range("D4").Value = 5.823916824656052
$ws.Range("E4").Value = 10.47746403678979
$ws.Range("G4").Value = 44.24503670219829
$ws.Range("H4").Value = 17.89627825539144
$ws.Range("I4").Value = 27.16094331229372
$ws.Range("L4").Value = 10.20911894682138
$ws.Range("M4").Value = 16.12089204758946
$ws.Range("N4").Value = 19.33186586191299
$ws.Range("B5").Value = 17.25348143649261
$ws.Range("C5").Value = 8.278312829030545
$ws.Range("D5").Value = 5.795957076022169
$ws.Range("E5").Value = 10.47725050203904
$ws.Range("G5").Value = 44.17672334632183
$ws.Range("H5").Value = 17.90212698487872
$ws.Range("I5").Value = 27.17635959772396
$ws.Range("L5").Value = 10.21266628178528
$ws.Range("M5").Value = 16.10294119095822
$ws.Range("N5").Value = 19.35123047021305
$ws.Range("B6").Value = 17.23526063896471
$ws.Range("C6").Value = 8.252598239377576
$ws.Range("D6").Value = 5.791330193725875
$ws.Range("E6").Value = 10.47721726347017
$ws.Range("G6").Value = 44.16561007958616
$ws.Range("H6").Value = 17.90316139923375
$ws.Range("I6").Value = 27.17901882646368
$ws.Range("L6").Value = 10.21327851079403
$ws.Range("M6").Value = 16.10000965690816
$ws.Range("N6").Value = 19.35447546170717
$ws.Range("B7").Value = 17.36174422498566
$ws.Range("C7").Value = 8.429306294610679
$ws.Range("D7").Value = 5.823538716320019
$ws.Range("E7").Value = 10.47746100821803
$ws.Range("G7").Value = 44.24410001241716
$ws.Range("H7").Value = 17.89635289231582
$ws.Range("I7").Value = 27.16114455618233
$ws.Range("L7").Value = 10.20916523238313
$ws.Range("M7").Value = 16.12064666681978
$ws.Range("N7").Value = 19.33212504238767
$ws.Range("B8").Value = 17.91894551834746
$ws.Range("C8").Value = 9.164854500492373
$ws.Range("D8").Value = 5.967676230329774
$ws.Range("E8").Value = 10.478859239737
$ws.Range("G8").Value = 44.62350902663126
$ws.Range("H8").Value = 17.87551885965781
$ws.Range("I8").Value = 27.09675756610119
$ws.Range("L8").Value = 10.19438545126653
$ws.Range("M8").Value = 16.21876801063176
$ws.Range("N8").Value = 19.23777597718578
$ws.Range("B9").Value = 18.99713470670969
$ws.Range("C9").Value = 10.44708828682585
$ws.Range("D9").Value = 6.254805301328537
$ws.Range("E9").Value = 10.48276265479756
$ws.Range("G9").Value = 45.48191910058652
$ws.Range("H9").Value = 17.86791279053964
$ws.Range("I9").Value = 27.02296210813679
$ws.Range("L9").Value = 10.17750972897146
$ws.Range("M9").Value = 16.43554919520026
$ws.Range("N9").Value = 19.06802217213474
$ws.Range("B10").Value = 19.77125921411616
$ws.Range("C10").Value = 11.29102080999303
$ws.Range("D10").Value = 6.466093588959263
$ws.Range("E10").Value = 10.4863412875044
$ws.Range("G10").Value = 46.17774429686991
$ws.Range("H10").Value = 17.88269284799925
$ws.Range("I10").Value = 27.0010511645252
$ws.Range("L10").Value = 10.17248026331658
$ws.Range("M10").Value = 16.60870048175372
$ws.Range("N10").Value = 18.95249961228227
$ws.Range("B11").Value = 20.11775384964022
$ws.Range("C11").Value = 11.65363273670293
$ws.Range("D11").Value = 6.561802981513459
$ws.Range("E11").Value = 10.48812428589902
$ws.Range("G11").Value = 46.50733497796281
$ws.Range("H11").Value = 17.89385016006397
$ws.Range("I11").Value = 26.99815505329829
$ws.Range("L11").Value = 10.17178831725013
$ws.Range("M11").Value = 16.69026029073899
$ws.Range("N11").Value = 18.90191853438709
$ws.Range("B12").Value = 20.24801543439036
$ws.Range("C12").Value = 11.78788462578293
$ws.Range("D12").Value = 6.597949005747368
$ws.Range("E12").Value = 10.48882179734938
$ws.Range("G12").Value = 46.63391953773405
$ws.Range("H12").Value = 17.89871262845035
$ws.Range("I12").Value = 26.99807874348592
$ws.Range("L12").Value = 10.17175528477869
$ws.Range("M12").Value = 16.72152620521
$ws.Range("N12").Value = 18.8830464300049
$ws.Range("B13").Value = 20.22000548787002
$ws.Range("C13").Value = 11.75910718997488
$ws.Range("D13").Value = 6.590169222467159
$ws.Range("E13").Value = 10.4886705822063
$ws.Range("G13").Value = 46.6065799748129
$ws.Range("H13").Value = 17.89763706672118
$ws.Range("I13").Value = 26.99804974700123
$ws.Range("L13").Value = 10.17175222342422
$ws.Range("M13").Value = 16.7147759208989
$ws.Range("N13").Value = 18.88709836422104
$ws.Range("B14").Value = 20.12849029127048
$ws.Range("C14").Value = 11.66473905942455
$ws.Range("D14").Value = 6.564778895707161
$ws.Range("E14").Value = 10.48818122290424
$ws.Range("G14").Value = 46.51771417490659
$ws.Range("H14").Value = 17.89423742582695
$ws.Range("I14").Value = 26.9981283073161
$ws.Range("L14").Value = 10.17178101351733
$ws.Range("M14").Value = 16.69282504089806
$ws.Range("N14").Value = 18.90036027502384
$ws.Range("B15").Value = 20.07230725697793
$ws.Range("C15").Value = 11.60653725042186
$ws.Range("D15").Value = 6.549212837257048
$ws.Range("E15").Value = 10.48788438499304
$ws.Range("G15").Value = 46.46350946869372
$ws.Range("H15").Value = 17.89223804722473
$ws.Range("I15").Value = 26.99830940565738
$ws.Range("L15").Value = 10.17182845330253
$ws.Range("M15").Value = 16.67942848725387
$ws.Range("N15").Value = 18.90852023839798
$ws.Range("B16").Value = 19.74848935903558
$ws.Range("C16").Value = 11.26689475032964
$ws.Range("D16").Value = 6.45982723183496
$ws.Range("E16").Value = 10.48622788722503
$ws.Range("G16").Value = 46.15645889352782
$ws.Range("H16").Value = 17.88205291137789
$ws.Range("I16").Value = 27.0013830426141
$ws.Range("L16").Value = 10.17255755744348
$ws.Range("M16").Value = 16.60342480450948
$ws.Range("N16").Value = 18.95584472911603
$ws.Range("B17").Value = 19.54828916202768
$ws.Range("C17").Value = 11.05308030405404
$ws.Range("D17").Value = 6.404859864019532
$ws.Range("E17").Value = 10.48525141613149
$ws.Range("G17").Value = 45.97136629671777
$ws.Range("H17").Value = 17.87694033634572
$ws.Range("I17").Value = 27.00508247703063
$ws.Range("L17").Value = 10.17341327490555
$ws.Range("M17").Value = 16.55749979706359
$ws.Range("N17").Value = 18.98538042378848
$ws.Range("B18").Value = 19.43261410531731
$ws.Range("C18").Value = 10.9280950388172
$ws.Range("D18").Value = 6.373208000253415
$ws.Range("E18").Value = 10.48470436890003
$ws.Range("G18").Value = 45.86614123644473
$ws.Range("H18").Value = 17.87441709067038
$ws.Range("I18").Value = 27.0078757070407
$ws.Range("L18").Value = 10.17405571834909
$ws.Range("M18").Value = 16.53134905920105
$ws.Range("N18").Value = 19.0025541582946
$ws.Range("B19").Value = 19.39336258944884
$ws.Range("C19").Value = 10.88543284963605
$ws.Range("D19").Value = 6.362486230520338
$ws.Range("E19").Value = 10.48452165355708
$ws.Range("G19").Value = 45.83072903531962
$ws.Range("H19").Value = 17.87363444227954
$ws.Range("I19").Value = 27.00893561314094
$ws.Range("L19").Value = 10.1742990595526
$ws.Range("M19").Value = 16.52254084565748
$ws.Range("N19").Value = 19.00840081082903
$ws.Range("B20").Value = 19.56965617961467
$ws.Range("C20").Value = 11.07604871112654
$ws.Range("D20").Value = 6.410715254381725
$ws.Range("E20").Value = 10.48535385256664
$ws.Range("G20").Value = 45.99094259448474
$ws.Range("H20").Value = 17.8774413816625
$ws.Range("I20").Value = 27.00461977063879
$ws.Range("L20").Value = 10.17330663467828
$ws.Range("M20").Value = 16.56236140487268
$ws.Range("N20").Value = 18.98221710266283
$ws.Range("B21").Value = 20.15539728510671
$ws.Range("C21").Value = 11.69254033010853
$ws.Range("D21").Value = 6.572239581789793
$ws.Range("E21").Value = 10.48832435334933
$ws.Range("G21").Value = 46.54376885276021
$ws.Range("H21").Value = 17.89521868672972
$ws.Range("I21").Value = 26.99807751439069
$ws.Range("L21").Value = 10.17176634685807
$ws.Range("M21").Value = 16.69926237368475
$ws.Range("N21").Value = 18.89645729314156
$ws.Range("B22").Value = 20.53262966012046
$ws.Range("C22").Value = 12.07761141442802
$ws.Range("D22").Value = 6.677222902924409
$ws.Range("E22").Value = 10.49039591085548
$ws.Range("G22").Value = 46.91537317853241
$ws.Range("H22").Value = 17.91055219953483
$ws.Range("I22").Value = 26.99975052097338
$ws.Range("L22").Value = 10.17209412596827
$ws.Range("M22").Value = 16.7909457718985
$ws.Range("N22").Value = 18.84205044423087
$ws.Range("B23").Value = 20.33184635895543
$ws.Range("C23").Value = 11.87372324573809
$ws.Range("D23").Value = 6.621256836081056
$ws.Range("E23").Value = 10.48927836115558
$ws.Range("G23").Value = 46.716133083639
$ws.Range("H23").Value = 17.90202867793579
$ws.Range("I23").Value = 26.99831229857227
$ws.Range("L23").Value = 10.17179727247219
$ws.Range("M23").Value = 16.74181729173452
$ws.Range("N23").Value = 18.87093865156349
$ws.Range("B24").Value = 19.55999793874415
$ws.Range("C24").Value = 11.06567111242443
$ws.Range("D24").Value = 6.408068186431567
$ws.Range("E24").Value = 10.48530749640712
$ws.Range("G24").Value = 45.98208844475768
$ws.Range("H24").Value = 17.87721356310912
$ws.Range("I24").Value = 27.00482688474455
$ws.Range("L24").Value = 10.17335437796472
$ws.Range("M24").Value = 16.56016268496889
$ws.Range("N24").Value = 18.98364663798229
$ws.Range("B25").Value = 18.70799009189648
$ws.Range("C25").Value = 10.11742631818002
$ws.Range("D25").Value = 6.176901372175259
$ws.Range("E25").Value = 10.48158205310049
$ws.Range("G25").Value = 45.23789769600648
$ws.Range("H25").Value = 17.86639866951317
$ws.Range("I25").Value = 27.03727326286002
$ws.Range("L25").Value = 10.18077978489348
$ws.Range("M25").Value = 16.3743922342711
$ws.Range("N25").Value = 19.11232232147839
